$wb = $excel.ActiveWorkbook

# --- 1. Clear the stray empty cell B2 on "ODI Batting" --------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("B2").ClearContents()

# --- 2. Add a new worksheet "ODI Batting Extra" after the last sheet ------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

# Reuse the same header style used by the other sheets (bold, bordered,
# centered) by copying formats from an existing header row.
$batting.Range("A1:F1").Copy()
$extra.Range("A1:F1").PasteSpecial(-4122)

# --- 3. Populate the header row --------------------------------------------
$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"

# --- 4. Populate the single data row ---------------------------------------
$extra.Range("A2").NumberFormat = "@"
$extra.Range("A2").Value = "4478"
$extra.Range("B2").Value = ""
$extra.Range("C2").Value = ""
$extra.Range("D2").Value = ""
$extra.Range("E2").Value = ""
$extra.Range("F2").Value = "NO"

Write-Output "edit applied"
